$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old placeholder rows (5, 6, 7) that held extra candidate names -
# the mentor records were consolidated down to rows 2-4.
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(5).Delete()

# Row 2 - Adam Jones
$ws.Range("E2").Value = "Male"
$ws.Range("H2").Value = "Bsc Chemistry, 2010, MSc BioChem, 2013"
$ws.Range("O2").Value = "No Preference"
$ws.Range("Q2").Value = "School of X"
$ws.Range("S2").Value = "Software Developer"
$ws.Range("T2").Value = "BAE Systems"
$ws.Range("U2").Value = "Tech"
$ws.Range("W2").Value = "Planning for the future and goal setting, Gaining insight to an industry/profession, Building a professional network, Writing/improving CVs, job applications and covering letters"

# Row 3 - Alice Sims
$ws.Range("E3").Value = "Female"
$ws.Range("H3").Value = "Bsc Computer Science, 2010, MSc Cyber Security, 2013"
$ws.Range("O3").Value = "Female"
$ws.Range("Q3").Value = "School of Y"
$ws.Range("R3").Value = "programming"
$ws.Range("S3").Value = "Software Developer"
$ws.Range("T3").Value = "Google"
$ws.Range("U3").Value = "IT"
$ws.Range("W3").Value = "Interview practice and preparation, Finding work experience (shadowing/internships/part-time work), Developing entrepreneurial skills, Support with setting up or growing a business"

# Row 4 - Mohammed Azar
$ws.Range("E4").Value = "Male"
$ws.Range("H4").Value = "Bsc Chemistry, 2010, MSc Business, 2013"
$ws.Range("O4").Value = "No Preference"
$ws.Range("Q4").Value = "School of S"
$ws.Range("R4").Value = "I do speak Urdu, interest in statistics."
$ws.Range("S4").Value = "Business Manager"
$ws.Range("T4").Value = "Lloyds Banking"
$ws.Range("U4").Value = "Banking"
$ws.Range("W4").Value = "Support with setting up or growing a business, Planning for the future and goal setting, Building a professional network, Interview practice and preparation, Finding work experience (shadowing/internships/part-time work), Developing entrepreneurial skills, Support with setting up or growing a business"

# Apply left/top alignment formatting to the new mentee-preference / career columns.
# (Multi-area ranges aren't reliable here, so format one cell directly, then
# copy/paste-special its format onto the remaining single-area ranges.)
$fmtCell = $ws.Range("O2")
$fmtCell.Font.Name = "Calibri"
$fmtCell.Font.Size = 11
$fmtCell.HorizontalAlignment = -4131
$fmtCell.VerticalAlignment = -4160

$fmtCell.Copy()
foreach ($addr in @("O3:O4", "Q2:Q4", "S2:S4", "U2:U4", "W2:W4")) {
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

$ws.Range("F4").Select()
